$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Change row 38 value from PAGAMENTO to SALÁRIO
$ws.Range("A38").Value = "SALÁRIO"

# Append new rows 204-207 with new data, copying the style of A203
$ws.Range("A204").Value = "SALÁRIO"
$ws.Range("A205").Value = "RESCISÃO TRABALHISTA"
$ws.Range("A206").Value = "PRESTADOR DE SERVIÇO"
$ws.Range("A207").Value = "PRESTADOR DE SERVIÇO"

$ws.Range("A203").Copy()
$ws.Range("A204:A207").PasteSpecial(-4122)
